{"js": "// Apply the \"midnight-madness\" content refresh:\n//  - New page title / H1 (also repeated later in the doc as a bold run)\n//  - Rewritten \"What we like\" / \"What we don't like\" bullet points\n//  - New meta-description style italic summary line\n// Replacements are done via Body.search(...) + Range.insertText(..., \"Replace\")\n// so formatting (bold/italic/list style) on each run is preserved.\n\nconst replacements = [\n  {\n    // Title appears twice verbatim: once as the Heading1 and once as a bold\n    // run near the end of the document. search() finds both occurrences.\n    find: \"Play Midnight Madness Free: Pros, Cons & Info | 2021 Review\",\n    replace: \"Play Midnight Madness Free Slot Game\",\n  },\n  {\n    find: \"Neon retro theme with flashy design\",\n    replace: \"Classic retro theme with neon graphics\",\n  },\n  {\n    find: \"Easy to understand gameplay for new players\",\n    replace: \"Mobile-friendly design\",\n  },\n  {\n    find: \"High variance with a solid RTP of 96%\",\n    replace: \"High variance for exciting gameplay\",\n  },\n  {\n    find: \"Mobile-friendly and available on desktop as well\",\n    replace: \"Potential maximum win of 15,000 times the bet\",\n  },\n  {\n    find: \"No bonus features or scatters available\",\n    replace: \"No bonus features or scatters\",\n  },\n  {\n    find: \"Maximum bet per slot is limited to 1.00 credit\",\n    replace: \"Limited maximum bet per slot\",\n  },\n  {\n    find: \"Play Midnight Madness for free: A high variance slot by Spearhead Studios with 96% RTP, no bonuses, and max win of 15k. Read our 2021 review for pros and cons.\",\n    replace: \"Read our review of Midnight Madness, a high variance slot game with a classic retro theme. Play for free.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"midnight-madness\" content refresh:\n#  - New page title / H1 (also repeated later in the doc as a bold run)\n#  - Rewritten \"What we like\" / \"What we don't like\" bullet points\n#  - New meta-description style italic summary line\n# Uses Find/Replace (wdReplaceAll) on the whole-document range so formatting\n# (bold/italic/list style) already on each run is preserved.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Midnight Madness Free: Pros, Cons & Info | 2021 Review\", \"Play Midnight Madness Free Slot Game\"),\n    @(\"Neon retro theme with flashy design\", \"Classic retro theme with neon graphics\"),\n    @(\"Easy to understand gameplay for new players\", \"Mobile-friendly design\"),\n    @(\"High variance with a solid RTP of 96%\", \"High variance for exciting gameplay\"),\n    @(\"Mobile-friendly and available on desktop as well\", \"Potential maximum win of 15,000 times the bet\"),\n    @(\"No bonus features or scatters available\", \"No bonus features or scatters\"),\n    @(\"Maximum bet per slot is limited to 1.00 credit\", \"Limited maximum bet per slot\"),\n    @(\"Play Midnight Madness for free: A high variance slot by Spearhead Studios with 96% RTP, no bonuses, and max win of 15k. Read our 2021 review for pros and cons.\", \"Read our review of Midnight Madness, a high variance slot game with a classic retro theme. Play for free.\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
